$wb = $excel.ActiveWorkbook

# --- Add the new "metadata" worksheet after the existing "data" sheet ---
$dataSheet = $wb.Worksheets.Item("data")
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# --- Header row (row 1) ---
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("B1:G1").Font.Bold = $true
$ws.Range("B1:G1").Borders.LineStyle = 1
$ws.Range("B1:G1").HorizontalAlignment = -4108
$ws.Range("B1:G1").VerticalAlignment = -4160

# --- Data row (row 2) ---
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

$ws.Range("B2").Value = "Congenital Diarrhoea"
$ws.Range("C2").Value = 89

# Write "1.7" as text (matches source workbook, which stores it as a string),
# then drop back to the default style so no stray number-format index sticks.
$ws.Range("D2").Value = "'1.7"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-08-15T02:24:56.013901Z"
$ws.Range("F2").Value = "2021-10-05 14:33:34.316419"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/89/?format=json"

# --- Refresh the "panel_query_time" (F column, "data" sheet) timestamps ---
$newTimes = @(
    "2021-10-05 14:33:34.320338",
    "2021-10-05 14:33:34.320348",
    "2021-10-05 14:33:34.320352",
    "2021-10-05 14:33:34.320354",
    "2021-10-05 14:33:34.320358",
    "2021-10-05 14:33:34.320360",
    "2021-10-05 14:33:34.320363",
    "2021-10-05 14:33:34.320366",
    "2021-10-05 14:33:34.320369",
    "2021-10-05 14:33:34.320372",
    "2021-10-05 14:33:34.320375",
    "2021-10-05 14:33:34.320378",
    "2021-10-05 14:33:34.320380",
    "2021-10-05 14:33:34.320383",
    "2021-10-05 14:33:34.320386",
    "2021-10-05 14:33:34.320389",
    "2021-10-05 14:33:34.320392",
    "2021-10-05 14:33:34.320395",
    "2021-10-05 14:33:34.320398",
    "2021-10-05 14:33:34.320400",
    "2021-10-05 14:33:34.320403",
    "2021-10-05 14:33:34.320406",
    "2021-10-05 14:33:34.320408",
    "2021-10-05 14:33:34.320411",
    "2021-10-05 14:33:34.320414",
    "2021-10-05 14:33:34.320416",
    "2021-10-05 14:33:34.320419",
    "2021-10-05 14:33:34.320422",
    "2021-10-05 14:33:34.320424",
    "2021-10-05 14:33:34.320427",
    "2021-10-05 14:33:34.320430",
    "2021-10-05 14:33:34.320433",
    "2021-10-05 14:33:34.320436",
    "2021-10-05 14:33:34.320438",
    "2021-10-05 14:33:34.320441",
    "2021-10-05 14:33:34.320444",
    "2021-10-05 14:33:34.320446",
    "2021-10-05 14:33:34.320449"
)

$data = $wb.Worksheets.Item("data")
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Keep "data" as the active tab (matches the workbook's original bookViews state).
$data.Activate()
